# Javascript_1.pptx edit:
#  - Slide 18 ("Undefined vs null"): add two clarifying lines to the
#    Content Placeholder:
#      * after the "Undefined" explanation paragraph -> ": default value of any variable" (bold)
#      * after the "Null" explanation paragraph       -> ": Something with no value" (bold)
#  - Handout master date field cache bumped by one day (8/18/2023 -> 8/19/2023)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 18: Undefined vs null
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(18)
$body = $slide.Shapes.Item(2)

# Locate a paragraph by (partial) text content. Re-scanning after every edit
# keeps this robust to paragraph-index shifts caused by earlier insertions.
function Find-ParagraphIndex($textFrameRange, $needle) {
    for ($i = 1; $i -le $textFrameRange.Paragraphs().Count; $i++) {
        if ($textFrameRange.Paragraphs($i).Text -like $needle) {
            return $i
        }
    }
    return 0
}

# --- "Undefined" paragraph: add "<colon><bold> default value of any variable </bold>"
$tr = $body.TextFrame.TextRange
$undefinedExplanationIdx = Find-ParagraphIndex $tr "*automatically assigned to a variable that has been declared*"
$undefPara = $tr.Paragraphs($undefinedExplanationIdx)
$undefPara.InsertAfter([char]13 + ": default value of any variable ") | Out-Null

# Re-fetch text range/paragraph after the mutation and style the new paragraph.
$tr = $body.TextFrame.TextRange
$newPara1 = $tr.Paragraphs($undefinedExplanationIdx + 1)
# Whole new paragraph starts out matching the previous paragraph's formatting;
# bold just the " default value of any variable " part, leave the leading ":" regular.
$boldRun1 = $newPara1.Characters(2, $newPara1.Length - 1)
$boldRun1.Font.Bold = $true

# --- "Null" paragraph: add bold "<colon> Something with no value "
$tr = $body.TextFrame.TextRange
$nullExplanationIdx = Find-ParagraphIndex $tr "*deliberate absence of any object value*"
$nullPara = $tr.Paragraphs($nullExplanationIdx)
$nullPara.InsertAfter([char]13 + ": Something with no value ") | Out-Null

$tr = $body.TextFrame.TextRange
$newPara2 = $tr.Paragraphs($nullExplanationIdx + 1)
$newPara2.Font.Bold = $true

# ---------------------------------------------------------------------------
# Handout master date placeholder: cached datetimeFigureOut field text.
# (Best effort -- real PowerPoint recalculates this automatically on save.)
# ---------------------------------------------------------------------------
try {
    $handout = $p.HandoutMaster
    $datePh = $handout.Shapes.Item(2)
    $datePh.TextFrame.TextRange.Text = "8/19/2023"
} catch {
}
